$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (bold font + thin border + centered alignment) from the
# existing last styled cell in column A down through the newly added rows,
# so the new rows match the style used by the rest of column A (style index 1).
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A12:A15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Write out the refreshed prediction table (rows re-sorted / re-sampled upstream).

# Row 2: even_MAG-GUT19787.fa
$ws.Range("A2").Value = "even_MAG-GUT19787.fa"
$ws.Range("B2").Value = [double]"0.9991853426428059"
$ws.Range("C2").Value = [double]"4.779999182348617e-07"
$ws.Range("D2").Value = [double]"0.0008074486677492747"
$ws.Range("E2").Value = [double]"2.220095147134157e-14"
$ws.Range("F2").Value = [double]"6.730689504325896e-06"
$ws.Range("G2").Value = [double]"0.9991853426428059"
$ws.Range("H2").Value = "f__Peptoniphilaceae"
$ws.Range("I2").Value = "f__Peptoniphilaceae"

# Row 3: even_MAG-GUT47840.fa
$ws.Range("A3").Value = "even_MAG-GUT47840.fa"
$ws.Range("B3").Value = [double]"0.9950079648239014"
$ws.Range("C3").Value = [double]"2.45729303589972e-05"
$ws.Range("D3").Value = [double]"0.003850736672401021"
$ws.Range("E3").Value = [double]"2.21999467764649e-14"
$ws.Range("F3").Value = [double]"0.001116725573316204"
$ws.Range("G3").Value = [double]"0.9950079648239014"
$ws.Range("H3").Value = "f__Peptoniphilaceae"
$ws.Range("I3").Value = "f__Peptoniphilaceae"

# Row 4: even_MAG-GUT49571.fa
$ws.Range("A4").Value = "even_MAG-GUT49571.fa"
$ws.Range("B4").Value = [double]"0.5190777206823617"
$ws.Range("C4").Value = [double]"0.2904083968279392"
$ws.Range("D4").Value = [double]"0.0008557517582532601"
$ws.Range("E4").Value = [double]"5.193396202267525e-13"
$ws.Range("F4").Value = [double]"0.1896581307309267"
$ws.Range("G4").Value = [double]"0.5190777206823617"
$ws.Range("H4").Value = "f__Peptoniphilaceae"
$ws.Range("I4").Value = "f__Peptoniphilaceae(reject)"

# Row 5: even_MAG-GUT5258.fa
$ws.Range("A5").Value = "even_MAG-GUT5258.fa"
$ws.Range("B5").Value = [double]"0.9962836617954696"
$ws.Range("C5").Value = [double]"2.476876221665816e-06"
$ws.Range("D5").Value = [double]"0.0006093971880639737"
$ws.Range("E5").Value = [double]"2.219952213447342e-14"
$ws.Range("F5").Value = [double]"0.00310446414022254"
$ws.Range("G5").Value = [double]"0.9962836617954696"
$ws.Range("H5").Value = "f__Peptoniphilaceae"
$ws.Range("I5").Value = "f__Peptoniphilaceae"

# Row 6: even_MAG-GUT7064.fa
$ws.Range("A6").Value = "even_MAG-GUT7064.fa"
$ws.Range("B6").Value = [double]"0.9979064276140123"
$ws.Range("C6").Value = [double]"8.073351499961667e-06"
$ws.Range("D6").Value = [double]"0.001722131690025905"
$ws.Range("E6").Value = [double]"2.219857605396516e-14"
$ws.Range("F6").Value = [double]"0.0003633673444395156"
$ws.Range("G6").Value = [double]"0.9979064276140123"
$ws.Range("H6").Value = "f__Peptoniphilaceae"
$ws.Range("I6").Value = "f__Peptoniphilaceae"

# Row 7: even_MAG-GUT7291.fa
$ws.Range("A7").Value = "even_MAG-GUT7291.fa"
$ws.Range("B7").Value = [double]"0.9956913986028839"
$ws.Range("C7").Value = [double]"8.239887700911545e-07"
$ws.Range("D7").Value = [double]"0.004177469353394843"
$ws.Range("E7").Value = [double]"2.219872620794984e-14"
$ws.Range("F7").Value = [double]"0.0001303080549290625"
$ws.Range("G7").Value = [double]"0.9956913986028839"
$ws.Range("H7").Value = "f__Peptoniphilaceae"
$ws.Range("I7").Value = "f__Peptoniphilaceae"

# Row 8: even_MAG-GUT88709.fa
$ws.Range("A8").Value = "even_MAG-GUT88709.fa"
$ws.Range("B8").Value = [double]"0.4741841439919613"
$ws.Range("C8").Value = [double]"2.73130476736386e-06"
$ws.Range("D8").Value = [double]"0.5257728302963622"
$ws.Range("E8").Value = [double]"2.28873607209728e-10"
$ws.Range("F8").Value = [double]"4.029417803531819e-05"
$ws.Range("G8").Value = [double]"0.5257728302963622"
$ws.Range("H8").Value = "f__Sporanaerobacteraceae"
$ws.Range("I8").Value = "f__Sporanaerobacteraceae(reject)"

# Row 9: even_MAG-GUT91256.fa
$ws.Range("A9").Value = "even_MAG-GUT91256.fa"
$ws.Range("B9").Value = [double]"0.9999973867881425"
$ws.Range("C9").Value = [double]"4.165575041322913e-07"
$ws.Range("D9").Value = [double]"1.754568398919589e-06"
$ws.Range("E9").Value = [double]"2.220445943010059e-14"
$ws.Range("F9").Value = [double]"4.420859322658218e-07"
$ws.Range("G9").Value = [double]"0.9999973867881425"
$ws.Range("H9").Value = "f__Peptoniphilaceae"
$ws.Range("I9").Value = "f__Peptoniphilaceae"

# Row 10: even_MAG-GUT91291.fa
$ws.Range("A10").Value = "even_MAG-GUT91291.fa"
$ws.Range("B10").Value = [double]"0.9990760963750309"
$ws.Range("C10").Value = [double]"5.922285073790967e-06"
$ws.Range("D10").Value = [double]"0.0009125964393434706"
$ws.Range("E10").Value = [double]"2.220122238265562e-14"
$ws.Range("F10").Value = [double]"5.38490052976206e-06"
$ws.Range("G10").Value = [double]"0.9990760963750309"
$ws.Range("H10").Value = "f__Peptoniphilaceae"
$ws.Range("I10").Value = "f__Peptoniphilaceae"

# Row 11: even_MAG-GUT91328.fa
$ws.Range("A11").Value = "even_MAG-GUT91328.fa"
$ws.Range("B11").Value = [double]"0.9964138112626825"
$ws.Range("C11").Value = [double]"1.154560353992336e-06"
$ws.Range("D11").Value = [double]"0.003448180833658645"
$ws.Range("E11").Value = [double]"2.219875564631368e-14"
$ws.Range("F11").Value = [double]"0.0001368533432825453"
$ws.Range("G11").Value = [double]"0.9964138112626825"
$ws.Range("H11").Value = "f__Peptoniphilaceae"
$ws.Range("I11").Value = "f__Peptoniphilaceae"

# Row 12: even_MAG-GUT91345.fa
$ws.Range("A12").Value = "even_MAG-GUT91345.fa"
$ws.Range("B12").Value = [double]"0.4257133207930784"
$ws.Range("C12").Value = [double]"0.2459000403154091"
$ws.Range("D12").Value = [double]"0.01272325994308441"
$ws.Range("E12").Value = [double]"1.207746760317836e-12"
$ws.Range("F12").Value = [double]"0.3156633789472204"
$ws.Range("G12").Value = [double]"0.4257133207930784"
$ws.Range("H12").Value = "f__Peptoniphilaceae"
$ws.Range("I12").Value = "f__Peptoniphilaceae(reject)"

# Row 13: even_MAG-GUT91672.fa
$ws.Range("A13").Value = "even_MAG-GUT91672.fa"
$ws.Range("B13").Value = [double]"0.9999981076734737"
$ws.Range("C13").Value = [double]"9.347136063404921e-08"
$ws.Range("D13").Value = [double]"1.667058982821181e-06"
$ws.Range("E13").Value = [double]"2.22044546734056e-14"
$ws.Range("F13").Value = [double]"1.317961607973203e-07"
$ws.Range("G13").Value = [double]"0.9999981076734737"
$ws.Range("H13").Value = "f__Peptoniphilaceae"
$ws.Range("I13").Value = "f__Peptoniphilaceae"

# Row 14: even_MAG-GUT91675.fa
$ws.Range("A14").Value = "even_MAG-GUT91675.fa"
$ws.Range("B14").Value = [double]"0.9999994315856907"
$ws.Range("C14").Value = [double]"8.294075811140223e-08"
$ws.Range("D14").Value = [double]"3.715505114134382e-07"
$ws.Range("E14").Value = [double]"2.220446008079375e-14"
$ws.Range("F14").Value = [double]"1.139230174586238e-07"
$ws.Range("G14").Value = [double]"0.9999994315856907"
$ws.Range("H14").Value = "f__Peptoniphilaceae"
$ws.Range("I14").Value = "f__Peptoniphilaceae"

# Row 15: even_MAG-GUT91866.fa
$ws.Range("A15").Value = "even_MAG-GUT91866.fa"
$ws.Range("B15").Value = [double]"0.06390333565106432"
$ws.Range("C15").Value = [double]"0.2989445139074004"
$ws.Range("D15").Value = [double]"0.01008423145424333"
$ws.Range("E15").Value = [double]"3.904225174202106e-13"
$ws.Range("F15").Value = [double]"0.6270679189869015"
$ws.Range("G15").Value = [double]"0.6270679189869015"
$ws.Range("H15").Value = "f__Tissierellaceae"
$ws.Range("I15").Value = "f__Tissierellaceae(reject)"
